$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 516, shifting existing rows 516-623 down to 518-625
$ws.Rows.Item(516).Resize(2).Insert()

# Populate new row 516
$ws.Cells.Item(516, 1).Value = 4
$ws.Cells.Item(516, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(516, 3).Value = "Los Lagos"
$ws.Cells.Item(516, 4).Value = 45244
$ws.Cells.Item(516, 5).Value = 10
$ws.Cells.Item(516, 6).Value = "Fruta"
$ws.Cells.Item(516, 7).Value = 100104
$ws.Cells.Item(516, 8).Value = "Frutos de pepita"
$ws.Cells.Item(516, 9).Value = 100104005
$ws.Cells.Item(516, 10).Value = "Pera"
$ws.Cells.Item(516, 11).Value = "Packham's Triumph"
$ws.Cells.Item(516, 12).Value = "Primera"
$ws.Cells.Item(516, 13).Value = 300
$ws.Cells.Item(516, 14).Value = 20000
$ws.Cells.Item(516, 15).Value = 20000
$ws.Cells.Item(516, 16).Value = 20000
$ws.Cells.Item(516, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(516, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(516, 19).Value = 1333
$ws.Cells.Item(516, 20).Value = 15

# Populate new row 517
$ws.Cells.Item(517, 1).Value = 4
$ws.Cells.Item(517, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(517, 3).Value = "Los Lagos"
$ws.Cells.Item(517, 4).Value = 45244
$ws.Cells.Item(517, 5).Value = 10
$ws.Cells.Item(517, 6).Value = "Fruta"
$ws.Cells.Item(517, 7).Value = 100104
$ws.Cells.Item(517, 8).Value = "Frutos de pepita"
$ws.Cells.Item(517, 9).Value = 100104005
$ws.Cells.Item(517, 10).Value = "Pera"
$ws.Cells.Item(517, 11).Value = "Packham's Triumph"
$ws.Cells.Item(517, 12).Value = "Segunda"
$ws.Cells.Item(517, 13).Value = 300
$ws.Cells.Item(517, 14).Value = 15000
$ws.Cells.Item(517, 15).Value = 15000
$ws.Cells.Item(517, 16).Value = 15000
$ws.Cells.Item(517, 17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(517, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(517, 19).Value = 1000
$ws.Cells.Item(517, 20).Value = 15
